$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add new row 54 with data: Date 08-Mar-2024 (serial 45359), 3 hours, "Navigation Session Variablen"
$ws.Cells.Item(54, 1).Value = 45359
$ws.Cells.Item(54, 1).NumberFormat = $ws.Cells.Item(53, 1).NumberFormat
$ws.Cells.Item(54, 2).Value = 3
$ws.Cells.Item(54, 3).Value = "Navigation Session Variablen"

# Move the selection to the newly added row, like a user would after typing the entry
$ws.Range("C54").Select()

$wb.Application.Calculate()
